$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = -0.206
$ws.Range("G2").Value = -0.186267845003399
$ws.Range("H2").Value = -0.1869476546566961
$ws.Range("I2").Value = -0.5152957171991842
$ws.Range("J2").Value = -0.5152957171991842
$ws.Range("K2").Value = -1.438
$ws.Range("L2").Value = -0.4887831407205983
$ws.Range("U2").Value = 0.06599999999999999
$ws.Range("V2").Value = 0.007440811724915445
$ws.Range("W2").Value = -0.8879588229177499
$ws.Range("X2").Value = 0.07104865067219132
$ws.Range("Y2").Value = -0.9590074735899412
$ws.Range("Z2").Value = 5.046312178387651
$ws.Range("AA2").Value = 0.01713756368689201
$ws.Range("AB2").Value = 0.06870152172747139
$ws.Range("AC2").Value = -0.05156395804057945
$ws.Range("AD2").Value = 0.521
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0.521
$ws.Range("AG2").Value = 0.455
$ws.Range("AH2").Value = 0.0554786497710574
$ws.Range("AI2").Value = -1.106157112526539
$ws.Range("AJ2").Value = 0.04879356568364612
$ws.Range("AK2").Value = -0.8472998137802609
$ws.Range("AL2").Value = 0.03
$ws.Range("AM2").Value = 0.02
$ws.Range("AN2").Value = -0.4846511627906976
$ws.Range("AO2").Value = -50.53333333333334
$ws.Range("AP2").Value = -0.4232558139534883
$ws.Range("AQ2").Value = -75.80000000000001

# Row 3
$ws.Range("D3").Value = -0.206
$ws.Range("G3").Value = -1.004132231404959
$ws.Range("H3").Value = -1.012396694214876
$ws.Range("I3").Value = -1.636363636363636
$ws.Range("J3").Value = -1.636363636363636
$ws.Range("K3").Value = -0.388
$ws.Range("L3").Value = -1.603305785123967
$ws.Range("V3").Value = -0.01282051282051282
$ws.Range("W3").Value = 1.190184049079755
$ws.Range("X3").Value = 0.07250438929973226
$ws.Range("Y3").Value = 1.117679659780022
$ws.Range("Z3").Value = -0.7908496732026143
$ws.Range("AA3").Value = 1.294117647058824
$ws.Range("AB3").Value = 0.06730946866450967
$ws.Range("AC3").Value = 1.226808178394314
$ws.Range("AD3").Value = 0.179
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0.179
$ws.Range("AG3").Value = 0.199
$ws.Range("AH3").Value = 0.1029327199539965
$ws.Range("AI3").Value = -0.3968957871396895
$ws.Range("AJ3").Value = 0.1131324616259238
$ws.Range("AK3").Value = -0.4617169373549883
$ws.Range("AL3").Value = 0.014
$ws.Range("AM3").Value = 0.014
$ws.Range("AN3").Value = -0.4578005115089513
$ws.Range("AO3").Value = -28.28571428571428
$ws.Range("AP3").Value = -0.5089514066496164
$ws.Range("AQ3").Value = -28.28571428571428

# Row 4
$ws.Range("G4").Value = -0.112962962962963
$ws.Range("H4").Value = -0.112962962962963
$ws.Range("I4").Value = -0.4148148148148148
$ws.Range("J4").Value = -0.4148148148148148
$ws.Range("K4").Value = -1.05
$ws.Range("L4").Value = -0.3888888888888889
$ws.Range("U4").Value = 0.08599999999999999
$ws.Range("V4").Value = 0.01176470588235294
$ws.Range("W4").Value = -2.966101694915254
$ws.Range("X4").Value = 0.06959291204465037
$ws.Range("Y4").Value = -3.035694606959905
$ws.Range("Z4").Value = 3.037120359955006
$ws.Range("AA4").Value = -1.25984251968504
$ws.Range("AB4").Value = 0.07009357479043311
$ws.Range("AC4").Value = -1.329936094475473
$ws.Range("AD4").Value = 0.342
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0.342
$ws.Range("AG4").Value = 0.256
$ws.Range("AH4").Value = 0.0446941975953999
$ws.Range("AI4").Value = -17.10000000000003
$ws.Range("AJ4").Value = 0.03383558022733281
$ws.Range("AK4").Value = -2.415094339622642
$ws.Range("AL4").Value = 0.016
$ws.Range("AM4").Value = 0.006
$ws.Range("AN4").Value = -0.5
$ws.Range("AO4").Value = -70
$ws.Range("AP4").Value = -0.3742690058479532
$ws.Range("AQ4").Value = -186.6666666666667
